$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The library-protocol code used for this run changed (E7760 -> E7420).
# Column K holds that value for every sample row (2-27); since it's the
# exact same text in every row it lives as one shared string, so this one
# write ripples to the whole column.
$ws.Range("K2:K27").Value = "E7420"

# Give the (now edited) protocol column its own distinct look rather than
# sharing the plain "Normal" style it inherited before.
$ws.Range("K2:K27").Font.Name = "Arial"
$ws.Range("K2:K27").Font.Size = 11

# roboticLibraryPrep was a hard-coded boolean; replace it with a live
# formula so it actually recalculates.
$ws.Range("L2:L27").Formula = "=FALSE()"

# Leave the selection on the column that was actually edited.
$null = $ws.Range("K2:K27").Select()
